$wb = $excel.ActiveWorkbook

# Add the new "missing" worksheet after "compact" (Worksheets.Add() inserts
# before the active sheet by default, so pass an explicit After: target).
$compact = $wb.Worksheets.Item("compact")
$missing = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $compact)
$missing.Name = "missing"

# Header row
$missing.Range("A1").Value = "Key"
$missing.Range("B1").Value = "Data.A"
$missing.Range("C1").Value = "AllNull"
$missing.Range("E1").Value = "Data.B"

# Row 2
$missing.Range("A2").Value = "SMITH"
$missing.Range("B2").Value = "Pull"
$missing.Range("E2").Value = 10

# Row 3
$missing.Range("A3").Value = "JOHNSON"
$missing.Range("B3").Value = "request"
$missing.Range("D3").Value = "should't read this"
$missing.Range("E3").Value = 15

# Row 4
$missing.Range("A4").Value = "NULLS"
$missing.Range("B4").Value = "issue"

# Row 5
$missing.Range("A5").Value = "MILLER"
$missing.Range("E5").Value = 35

# Row 6 intentionally left blank

# Row 7
$missing.Range("A7").Value = "MICHEAL"
$missing.Range("B7").Value = "after"
$missing.Range("E7").Value = 50

# Selection bookkeeping to mirror the authored file: the "compact" sheet
# keeps a plain (non-active) selection, and "missing" becomes the
# selected/active tab with its own selection.
[void]$compact.Range("A2:A10").Select()
[void]$missing.Select()
[void]$missing.Range("E11").Select()
